$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated optimization results reflecting the new 2% return target.
# Column C = "Opt Portfolio", Column D = "Opt Portfolio with View"

$ws.Range("C2").Value = 0.000000000000000003876022766635678
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = 0.02124815174069435
$ws.Range("D3").Value = 0.02124764584718219

$ws.Range("C4").Value = 0.1273388301359862
$ws.Range("D4").Value = 0.1273383454430054

$ws.Range("C5").Value = 0.1697700754744735
$ws.Range("D5").Value = 0.1697703529919921

$ws.Range("C6").Value = 0.2390017860881791
$ws.Range("D6").Value = 0.2390020959213678

$ws.Range("C7").Value = 0.09744088543812593
$ws.Range("D7").Value = 0.09744157541012331

$ws.Range("C8").Value = 0.345200271122541
$ws.Range("D8").Value = 0.3451999843863293
